$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Promotion P001)
$ws.Range("C2").Value = "L'art d'évoluer à sa source"
$ws.Range("D2").Value = "Montant Fixe"
$ws.Range("E2").Value = 45559
$ws.Range("F2").Value = 45606

# Row 3 (Promotion P002)
$ws.Range("C3").Value = "L'art de concrétiser vos projets plus facilement"
$ws.Range("D3").Value = "Réduction"
$ws.Range("E3").Value = 45358
$ws.Range("F3").Value = 45545

# Row 4 (Promotion P003)
$ws.Range("C4").Value = "L'avantage de changer de manière efficace"
$ws.Range("E4").Value = 45518
$ws.Range("F4").Value = 45551

# Row 5 (Promotion P004)
$ws.Range("C5").Value = "L'art de concrétiser vos projets autrement"
$ws.Range("D5").Value = "Produit Offert"
$ws.Range("E5").Value = 45307
$ws.Range("F5").Value = 45481

# Row 6 (Promotion P005)
$ws.Range("C6").Value = "Le pouvoir d'évoluer de manière efficace"
$ws.Range("D6").Value = "Montant Fixe"
$ws.Range("E6").Value = 45539
$ws.Range("F6").Value = 45554

# Row 7 (Promotion P006)
$ws.Range("C7").Value = "Le pouvoir d'évoluer autrement"
$ws.Range("E7").Value = 45410
$ws.Range("F7").Value = 45474

# Row 8 (Promotion P007)
$ws.Range("C8").Value = "L'art d'avancer de manière sûre"
$ws.Range("E8").Value = 45364
$ws.Range("F8").Value = 45400

# Row 9 (Promotion P008)
$ws.Range("C9").Value = "Le plaisir d'avancer en toute tranquilité"
$ws.Range("E9").Value = 45288
$ws.Range("F9").Value = 45552

# Row 10 (Promotion P009)
$ws.Range("C10").Value = "La liberté de concrétiser vos projets en toute tranquilité"
$ws.Range("D10").Value = "Réduction"
$ws.Range("E10").Value = 45618
$ws.Range("F10").Value = 45625

# Row 11 (Promotion P010)
$ws.Range("C11").Value = "La liberté d'innover à sa source"
$ws.Range("D11").Value = "Réduction"
$ws.Range("E11").Value = 45473
$ws.Range("F11").Value = 45525
